# Auto-generated Excel COM-interop script
# Applies the "Updated symbol list on Thu Dec 29 17:30:17 UTC 2022 with GitHub Actions" commit:
# refreshed Price (column D) quotes and a handful of Best/Worst-in-24h badges
# appended/removed from the ranking label (column E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Price 246.49 -> 246.14
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.14"

# Row 3: Price 24.20 -> 24.21
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.21"

# Row 4: Price 5.294 -> 5.298
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.298"

# Row 5: Price 0.05794 -> 0.05793
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05793"

# Row 6: Price 6.478 -> 6.476
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.476"

# Row 7: Price 3.136 -> 3.138
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.138"

# Row 8: Price 0.8168 -> 0.8175
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8175"

# Row 9: Price 0.8771 -> 0.8697
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8697"

# Row 11: Price 0.07004 -> 0.07008
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07008"

# Row 12: Price 0.03134 -> 0.03135
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03135"

# Row 13: Price 0.02917 -> 0.02918
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02918"

# Row 14: Price 0.09415 -> 0.09402
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09402"

# Row 15: Price 3.748 -> 3.741
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.741"

# Row 16: Price 0.001529 -> 0.001531
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001531"

# Row 17: Price 0.04668 -> 0.04687
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04687"

# Row 18: Price 0.0005995 -> 0.0005975
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005975"
# Row 18: ranking label '17OneONEWorstin24h' -> '17OneONE'
$ws.Range("E18").Value = "17OneONE"

# Row 19: Price 0.006052 -> 0.006102
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006102"

# Row 20: Price 0.001245 -> 0.001244
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001244"

# Row 21: Price 0.004664 -> 0.004663
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004663"

# Row 22: Price 0.00006102 -> 0.00006106
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006106"
# Row 22: ranking label '21NitroExNTX' -> '21NitroExNTXWorstin24h'
$ws.Range("E22").Value = "21NitroExNTXWorstin24h"

# Row 23: Price 3.533 -> 3.530
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.530"

# Row 24: Price 2.144 -> 2.152
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.152"

# Row 25: Price 0.3189 -> 0.3188
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3188"

# Row 26: Price 0.1306 -> 0.1305
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1305"

# Row 40: Price 0.03723 -> 0.03725
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03725"

# Row 41: Price 0.006372 -> 0.006436
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006436"
# Row 41: ranking label '40KickTokenKICK' -> '40KickTokenKICKBestin24h'
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# Row 42: Price 0.1057 -> 0.1056
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1056"

# Row 43: Price 0.003001 -> 0.003003
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003003"
# Row 43: ranking label '42CEJICEJIBestin24h' -> '42CEJICEJI'
$ws.Range("E43").Value = "42CEJICEJI"

# Row 45: Price 0.00005275 -> 0.00005257
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005257"

# Row 46: Price 0.00000000750 -> 0.00000000751
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"

# Row 49: Price 0.00002101 -> 0.00002102
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"

# Row 50: Price 0.0002001 -> 0.0002002
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
